$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: A1:E1 updated (F1 and beyond unchanged)
$ws.Range("A1").Value = 0.00084860496271119048
$ws.Range("B1").Value = 0.00084860496271122962
$ws.Range("C1").Value = 0.99915139503728889
$ws.Range("D1").Value = 0.99915139503728889
$ws.Range("E1").Value = 0.00084860496271118668

# Row 2: A2:E2 updated (F2 and beyond unchanged)
$ws.Range("A2").Value = 0.99915139503728889
$ws.Range("B2").Value = 0.99915139503728889
$ws.Range("C2").Value = 0.00084860496271118668
$ws.Range("D2").Value = 0.00084860496271118668
$ws.Range("E2").Value = 0.99915139503728889
